# CIERRE 5 OCT 2021
# Fill in the September-credits ledger ("Hoja1") rows 14-19 with the
# remittances that cleared between 26-Sep-2021 and 01-Oct-2021, and leave
# the selection where the user's cursor ended up (G20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 14 - MAURO
$ws.Range("A14").Value = 44465
$ws.Range("D14").Value = "MAURO"
$ws.Range("E14").Value = 3670
$ws.Range("F14").Value = 44467
$ws.Range("G14").Value = 3670

# Row 15 - GUSTAVO
$ws.Range("A15").Value = 44465
$ws.Range("D15").Value = "GUSTAVO"
$ws.Range("E15").Value = 7851
$ws.Range("F15").Value = 44466
$ws.Range("G15").Value = 7851

# Row 16 - GUSTAVO
$ws.Range("A16").Value = 44466
$ws.Range("D16").Value = "GUSTAVO"
$ws.Range("E16").Value = 4662
$ws.Range("F16").Value = 44467
$ws.Range("G16").Value = 4662

# Row 17 - GUSTAVO
$ws.Range("A17").Value = 44467
$ws.Range("D17").Value = "GUSTAVO"
$ws.Range("E17").Value = 2213
$ws.Range("F17").Value = 44469
$ws.Range("G17").Value = 2213

# Row 18 - GUSTAVO
$ws.Range("A18").Value = 44469
$ws.Range("D18").Value = "GUSTAVO"
$ws.Range("E18").Value = 6855
$ws.Range("F18").Value = 44470
$ws.Range("G18").Value = 6855

# Row 19 - EL PRIMO
$ws.Range("A19").Value = 44470
$ws.Range("D19").Value = "EL PRIMO"
$ws.Range("E19").Value = 14800
$ws.Range("F19").Value = 44470
$ws.Range("G19").Value = 14800

# Leave the selection on G20, matching where the user finished editing.
$ws.Range("G20").Select() | Out-Null
